$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Time Record")

# Fill in the Thursday/Friday/Saturday/Sunday clock-in/clock-out times
# (time-of-day serials: fraction of a 24h day, matching the sheet's existing style)
$ws.Range("C16").Value = 13/24
$ws.Range("D16").Value = 16/24

$ws.Range("C17").Value = 10/24
$ws.Range("D17").Value = 11/24

$ws.Range("C18").Value = 20/24
$ws.Range("D18").Value = 23/24

$ws.Range("C19").Value = 20/24
$ws.Range("D19").Value = 23/24

# Adjust column B width slightly narrower
$ws.Columns.Item(2).ColumnWidth = 11.3

# Re-set the print area (adds another Print_Area defined name entry, matching Excel's append-only behavior)
$ws.PageSetup.PrintArea = '$A$1:$K$27'

# Update the active selection to match
$ws.Range("M17").Select()
